# Applies updated values to columns M and N ("Correct Extracted Predicates with
# Parents" and "Correct Extracted Predicates with Related") for several rows,
# reflecting corrected prediction extraction results (Gen1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = @{ M = 2; N = 2 }
    10 = @{ N = 3 }
    19 = @{ M = 1; N = 1 }
    30 = @{ M = 0 }
    31 = @{ M = 2; N = 2 }
    32 = @{ M = 1 }
    33 = @{ M = 3; N = 3 }
    42 = @{ M = 4; N = 4 }
    45 = @{ M = 2; N = 2 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
